$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 2512.9167
$ws.Range("I100").Value = 2087.8333
$ws.Range("J100").Value = 2938
$ws.Range("K100").Value = 2087.8333
$ws.Range("L100").Value = 2938
$ws.Range("M100").Value = -1546.8333
$ws.Range("N100").Value = -4020
$ws.Range("H107").Value = 802.875
$ws.Range("I107").Value = 219.66667
$ws.Range("J107").Value = 2552.5
$ws.Range("K107").Value = 219.66667
$ws.Range("L107").Value = 2552.5
$ws.Range("M107").Value = 1700.33333
$ws.Range("N107").Value = -6392.5
$ws.Range("H111").Value = 790.375
$ws.Range("I111").Value = 389.8
$ws.Range("J111").Value = 1458
$ws.Range("K111").Value = 1169.4
$ws.Range("L111").Value = 4374
$ws.Range("M111").Value = 1897.6
$ws.Range("N111").Value = -10508
$ws.Range("H132").Value = 1785.6111
$ws.Range("I132").Value = 1299.5483
$ws.Range("K132").Value = 3898.6449
$ws.Range("M132").Value = -1368.6449
$ws.Range("H133").Value = 88213.75
$ws.Range("J133").Value = 88213.75
$ws.Range("L133").Value = 88213.75
$ws.Range("N133").Value = -98333.75
$ws.Range("H134").Value = 52493.75
$ws.Range("J134").Value = 55421.43
$ws.Range("L134").Value = 55421.43
$ws.Range("N134").Value = -65561.42999999999
$ws.Range("H135").Value = 1257.8334
$ws.Range("I135").Value = 1303.2
$ws.Range("K135").Value = 11728.8
$ws.Range("M135").Value = -9193.800000000001
$ws.Range("H137").Value = 1117088.8
$ws.Range("I137").Value = 1689.2222
$ws.Range("K137").Value = 5067.6666
$ws.Range("M137").Value = -2517.6666

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H105").Value = 117522.664
$ws.Range("J105").Value = 117522.664
$ws.Range("L105").Value = 117522.664
$ws.Range("N105").Value = -124510.664

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 49563.43
$ws.Range("I20").Value = 68256.87
$ws.Range("K20").Value = 68256.87
$ws.Range("M20").Value = -68009.87
$ws.Range("H95").Value = 15087
$ws.Range("J95").Value = 14504.2
$ws.Range("L95").Value = 14504.2
$ws.Range("N95").Value = -19996.2
$ws.Range("H99").Value = 3236077
$ws.Range("I99").Value = 151072.58
$ws.Range("J99").Value = 10434420
$ws.Range("K99").Value = 151072.58
$ws.Range("L99").Value = 10434420
$ws.Range("M99").Value = -149574.58
$ws.Range("N99").Value = -10437416
$ws.Range("H103").Value = 15532.667
$ws.Range("J103").Value = 15532.667
$ws.Range("L103").Value = 15532.667
$ws.Range("N103").Value = -17876.667

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H14").Value = 16474.75
$ws.Range("I14").Value = 16474.75
$ws.Range("K14").Value = 16474.75
$ws.Range("M14").Value = -16304.75
$ws.Range("H86").Value = 250002130
$ws.Range("I86").Value = 250002130
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 250002130
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -250001007
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 250002130
$ws.Range("I89").Value = 250002130
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 1250010650
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -1250005034
$ws.Range("N89").ClearContents()
$ws.Range("H141").Value = 204781.83
$ws.Range("I141").Value = 60000
$ws.Range("J141").Value = 213298.4
$ws.Range("K141").Value = 60000
$ws.Range("L141").Value = 213298.4
$ws.Range("N141").Value = -223658.4
$ws.Range("M141").Value = -54820

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1682.6666
$ws.Range("I5").Value = 1428.8
$ws.Range("K5").Value = 4286.4
$ws.Range("M5").Value = -4174.4
$ws.Range("H80").Value = 6749.75
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 6749.75
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 20249.25
$ws.Range("M80").ClearContents()
$ws.Range("N80").Value = -22121.25
$ws.Range("H83").Value = 6749.75
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 6749.75
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 60747.75
$ws.Range("M83").ClearContents()
$ws.Range("N83").Value = -70107.75
$ws.Range("H121").Value = 1642.9722
$ws.Range("I121").Value = 1387.7142
$ws.Range("J121").Value = 1704.5862
$ws.Range("K121").Value = 4163.142599999999
$ws.Range("L121").Value = 5113.7586
$ws.Range("M121").Value = -2853.142599999999
$ws.Range("N121").Value = -7733.7586
$ws.Range("H129").Value = 4856.75
$ws.Range("J129").Value = 5966
$ws.Range("L129").Value = 17898
$ws.Range("N129").Value = -27898
$ws.Range("H131").Value = 37372.855
$ws.Range("J131").Value = 1780.44
$ws.Range("L131").Value = 5341.32
$ws.Range("N131").Value = -15421.32
$ws.Range("H135").Value = 1682.6666
$ws.Range("I135").Value = 1428.8
$ws.Range("K135").Value = 12859.2
$ws.Range("M135").Value = -10324.2

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H95").Value = 30999.2
$ws.Range("J95").Value = 30999.2
$ws.Range("L95").Value = 30999.2
$ws.Range("N95").Value = -36491.2
$ws.Range("H102").Value = 796.7406999999999
$ws.Range("I102").Value = 596.4583
$ws.Range("K102").Value = 596.4583
$ws.Range("M102").Value = 1025.5417
$ws.Range("H113").Value = 2565759.8
$ws.Range("J113").Value = 4764051
$ws.Range("L113").Value = 4764051
$ws.Range("N113").Value = -4768391
$ws.Range("H132").Value = 4360.3667
$ws.Range("I132").Value = 3932.3635
$ws.Range("K132").Value = 11797.0905
$ws.Range("M132").Value = -9267.0905
$ws.Range("H136").Value = 36751.684
$ws.Range("J136").Value = 36751.684
$ws.Range("L136").Value = 110255.052
$ws.Range("N136").Value = -115355.052

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 3799.5
$ws.Range("I16").Value = 4374.5
$ws.Range("K16").Value = 4374.5
$ws.Range("M16").Value = -4204.5
$ws.Range("H93").Value = 2079.6667
$ws.Range("I93").Value = 2084.353
$ws.Range("K93").Value = 2084.353
$ws.Range("M93").Value = -836.3530000000001
$ws.Range("H106").Value = 19685
$ws.Range("J106").Value = 19685
$ws.Range("L106").Value = 19685
$ws.Range("N106").Value = -22209
$ws.Range("H132").Value = 16640.7
$ws.Range("I132").Value = 35714
$ws.Range("J132").Value = 3925.1667
$ws.Range("K132").Value = 107142
$ws.Range("L132").Value = 11775.5001
$ws.Range("M132").Value = -104612
$ws.Range("N132").Value = -16835.5001
$ws.Range("H136").Value = 4912.1924
$ws.Range("I136").Value = 5704.4614
$ws.Range("J136").Value = 4119.923
$ws.Range("K136").Value = 17113.3842
$ws.Range("L136").Value = 12359.769
$ws.Range("M136").Value = -14563.3842
$ws.Range("N136").Value = -17459.769

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H12").Value = 3999.5
$ws.Range("I12").Value = 3999
$ws.Range("J12").Value = 4000
$ws.Range("K12").Value = 3999
$ws.Range("L12").Value = 4000
$ws.Range("M12").Value = -3857
$ws.Range("N12").Value = -4284
$ws.Range("H46").Value = 153166.44
$ws.Range("J46").Value = 153166.44
$ws.Range("L46").Value = 153166.44
$ws.Range("N46").Value = -153628.44
$ws.Range("H62").Value = 30000
$ws.Range("J62").Value = 30000
$ws.Range("L62").Value = 30000
$ws.Range("N62").Value = -31248
$ws.Range("H65").Value = 30000
$ws.Range("J65").Value = 30000
$ws.Range("L65").Value = 150000
$ws.Range("N65").Value = -156240
$ws.Range("H94").Value = 23938.3
$ws.Range("I94").Value = 25694
$ws.Range("J94").Value = 23499.375
$ws.Range("K94").Value = 25694
$ws.Range("L94").Value = 23499.375
$ws.Range("M94").Value = -24793
$ws.Range("N94").Value = -25301.375
$ws.Range("H107").Value = 1781.8723
$ws.Range("I107").Value = 1250
$ws.Range("J107").Value = 2211.4614
$ws.Range("K107").Value = 3750
$ws.Range("L107").Value = 6634.3842
$ws.Range("M107").Value = -1830
$ws.Range("N107").Value = -10474.3842
$ws.Range("H132").Value = 2719068
$ws.Range("I132").Value = 1303
$ws.Range("J132").Value = 7248676.5
$ws.Range("K132").Value = 3909
$ws.Range("L132").Value = 21746029.5
$ws.Range("M132").Value = -1379
$ws.Range("N132").Value = -21751089.5
$ws.Range("H134").Value = 153166.44
$ws.Range("J134").Value = 153166.44
$ws.Range("L134").Value = 459499.32
$ws.Range("N134").Value = -464569.32
